$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.810.78"
$ws.Range("E2").Value = "  +5.53%  "
$ws.Range("D3").Value = "2.424.85"
$ws.Range("E3").Value = "  +5.70%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'563.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.48%  "
$ws.Range("D6").Value = "'140.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +9.18%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  +2.94%  "
$ws.Range("D9").Value = "2.423.40"
$ws.Range("E9").Value = "  +5.83%  "
$ws.Range("D10").Value = "'0.104"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.45%  "
$ws.Range("E11").Value = "  +3.82%  "
$ws.Range("D13").Value = "'0.349"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.28%  "
$ws.Range("D14").Value = "'26.26"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +13.93%  "
$ws.Range("D15").Value = "2.860.27"
$ws.Range("E15").Value = "  +5.65%  "
$ws.Range("D16").Value = "62.700.95"
$ws.Range("E16").Value = "  +5.40%  "
$ws.Range("E17").Value = "  +8.13%  "
$ws.Range("D18").Value = "2.430.01"
$ws.Range("E18").Value = "  +5.75%  "
$ws.Range("D19").Value = "'11.19"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.84%  "
$ws.Range("D20").Value = "'339.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +10.07%  "
$ws.Range("E21").Value = "  +4.64%  "
$ws.Range("D22").Value = "'6.78"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.58%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").Value = "'5.66"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("D25").Value = "'65.18"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.61%  "
$ws.Range("E26").Value = "  +3.42%  "
$ws.Range("E27").Value = "  -0.25%  "
$ws.Range("E28").Value = "  +14.98%  "
$ws.Range("D29").Value = "'8.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.38%  "
$ws.Range("E30").Value = "  +13.10%  "
$ws.Range("D31").Value = "0.0₃0793"
$ws.Range("E31").Value = "  +10.55%  "
$ws.Range("E32").Value = "  +7.23%  "
$ws.Range("D33").Value = "'6.50"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +13.00%  "
$ws.Range("D34").Value = "'173.77"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.17%  "
$ws.Range("E35").Value = "  +10.31%  "
$ws.Range("E36").Value = "  +5.72%  "
$ws.Range("D37").Value = "'18.57"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.88%  "
$ws.Range("D38").Value = "'370.50"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +19.22%  "
$ws.Range("E39").Value = "  +12.68%  "
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("D41").Value = "'0.997"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.40%  "
$ws.Range("D42").Value = "'1.69"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +13.57%  "
$ws.Range("D43").Value = "'39.83"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.59%  "
$ws.Range("D44").Value = "'145.82"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'3.68"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.38%  "
$ws.Range("D46").Value = "'20.50"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +11.16%  "
$ws.Range("D47").Value = "'0.589"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.97%  "
$ws.Range("E48").Value = "  +1.45%  "
$ws.Range("E49").Value = "  +6.73%  "
$ws.Range("E50").Value = "  +5.91%  "
$ws.Range("D51").Value = "'17.81"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.14%  "
